# Update odds data on Sheet1 to reflect latest FlashScore values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("R2").Value = 1.63

# Row 5 updates
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("Q5").Value = 1.77
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.83
$ws.Range("AD5").Value = 7.5
$ws.Range("AE5").Value = 17
$ws.Range("AN5").Value = 7.5
$ws.Range("AR5").Value = 126
